{"js": "// Remove the trailing \"Ver no Jupiter ...\" / \"\u00a9 2020 ...\" footer block,\n// along with the blank paragraph that precedes it, from the end of the\n// document (right after the \"LOQ4038: ...\" requisito paragraph).\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items/text\");\nawait context.sync();\n\nconst items = paragraphs.items;\n\n// Find the \"Ver no Jupiter ...\" paragraph and, from there, the \"\u00a9 2020 ...\"\n// (Jekyll) paragraph that follows it. Only act if both distinctive\n// paragraphs are actually present, so the script is a safe no-op if the\n// footer block has already been removed.\nlet startIdx = -1;\nlet endIdx = -1;\nfor (let i = 0; i < items.length; i++) {\n  if (startIdx === -1 && items[i].text.indexOf(\"Ver no Jupiter\") !== -1) {\n    startIdx = i;\n  }\n  if (startIdx !== -1 && items[i].text.indexOf(\"Powered by Jekyll\") !== -1) {\n    endIdx = i;\n    break;\n  }\n}\n\nif (startIdx !== -1 && endIdx !== -1) {\n  const toDelete = [];\n  // Include the blank spacer paragraph immediately before \"Ver no Jupiter ...\"\n  if (startIdx - 1 >= 0 && items[startIdx - 1].text.trim() === \"\") {\n    toDelete.push(items[startIdx - 1]);\n  }\n  for (let i = startIdx; i <= endIdx; i++) {\n    toDelete.push(items[i]);\n  }\n\n  for (const p of toDelete) {\n    p.delete();\n  }\n\n  await context.sync();\n}\n", "ps1": "# Remove the trailing \"Ver no Jupiter ...\" / \"\u00a9 2020 ...\" footer block,\n# along with the blank paragraph that precedes it, from the end of the\n# document (right after the \"LOQ4038: ...\" requisito paragraph).\n$d = $word.ActiveDocument\n\n# Find the first paragraph that contains \"Ver no Jupiter\" and the first one\n# (from there on) that contains \"Powered by Jekyll\" (the \"\u00a9 2020 ...\" line).\n$startPara = $null\n$endPara = $null\nforeach ($p in $d.Paragraphs) {\n    $t = $p.Range.Text\n    if ($startPara -eq $null -and $t -match \"Ver no Jupiter\") {\n        $startPara = $p\n    }\n    if ($startPara -ne $null -and $t -match \"Powered by Jekyll\") {\n        $endPara = $p\n        break\n    }\n}\n\nif ($startPara -ne $null -and $endPara -ne $null) {\n    # Include the blank spacer paragraph right before \"Ver no Jupiter ...\" too.\n    $deleteStart = $startPara.Range.Start\n    $prev = $startPara.Previous()\n    if ($prev -ne $null -and $prev.Range.Text.Trim() -eq \"\") {\n        $deleteStart = $prev.Range.Start\n    }\n\n    # Build one contiguous range covering every paragraph mark that must go\n    # (deleting paragraph-by-paragraph can strand an extra pilcrow behind).\n    $rng = $d.Range($deleteStart, $endPara.Range.End)\n    $rng.Delete()\n}\n"}
